$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("parameters")

# --- Update TV/LV/HV (D, E, F) input values for rows 2-28 ---
$ws.Range("D2").Value = 16.801971423166499
$ws.Range("E2").Value = 0.15836293333333301
$ws.Range("F2").Value = 3.9590733333333299
$ws.Range("D3").Value = 0.001
$ws.Range("E3").Value = 0.001
$ws.Range("F3").Value = 100
$ws.Range("D4").Value = 0.1
$ws.Range("E4").Value = 0.1
$ws.Range("F4").Value = 2
$ws.Range("D5").Value = 929.80869219799104
$ws.Range("E5").Value = 0.25
$ws.Range("F5").Value = 4
$ws.Range("D6").Value = 0.57342071835498398
$ws.Range("E6").Value = 0.25
$ws.Range("F6").Value = 4
$ws.Range("D7").Value = 14463.6907675243
$ws.Range("E7").Value = 0.25
$ws.Range("F7").Value = 4
$ws.Range("D8").Value = 0.27914923181321899
$ws.Range("E8").Value = 0.25
$ws.Range("F8").Value = 4
$ws.Range("D9").Value = 129140.09613860901
$ws.Range("E9").Value = 0.25
$ws.Range("F9").Value = 4
$ws.Range("D10").Value = 1667.2259802286201
$ws.Range("E10").Value = 0.25
$ws.Range("F10").Value = 4
$ws.Range("D11").Value = 33.626295339999999
$ws.Range("E11").Value = 0.25
$ws.Range("F11").Value = 4
$ws.Range("D12").Value = 26.672409326424798
$ws.Range("E12").Value = 0.25
$ws.Range("F12").Value = 4
$ws.Range("D13").Value = 97.625215889464499
$ws.Range("E13").Value = 0.47499999999999998
$ws.Range("F13").Value = 1.73035714285714
$ws.Range("D14").Value = 9.6654488662369307
$ws.Range("E14").Value = 0.1
$ws.Range("F14").Value = 8
$ws.Range("D15").Value = 171.69402414356401
$ws.Range("E15").Value = 0.5
$ws.Range("F15").Value = 4
$ws.Range("D16").Value = 4.5149999999999997
$ws.Range("E16").Value = 0.56855746651966299
$ws.Range("F16").Value = 1.43144253348033
$ws.Range("D17").Value = 0.49622999999999901
$ws.Range("E17").Value = 0.52772808586761999
$ws.Range("F17").Value = 1.4722719141323699
$ws.Range("D18").Value = 0.15879359999999901
$ws.Range("E18").Value = 0.52772808586761999
$ws.Range("F18").Value = 1.4722719141323699
$ws.Range("D19").Value = 1
$ws.Range("E19").Value = 1
$ws.Range("F19").Value = 1
$ws.Range("D20").Value = 1
$ws.Range("E20").Value = 1
$ws.Range("F20").Value = 1
$ws.Range("D21").Value = 1
$ws.Range("E21").Value = 1
$ws.Range("F21").Value = 1
$ws.Range("D22").Value = 1
$ws.Range("E22").Value = 1
$ws.Range("F22").Value = 1
$ws.Range("D23").Value = 1
$ws.Range("E23").Value = 1
$ws.Range("F23").Value = 1
$ws.Range("D24").Value = 0.15
$ws.Range("E24").Value = 0.25
$ws.Range("F24").Value = 39.461538461538403
$ws.Range("D25").Value = 58.290155440414502
$ws.Range("E25").Value = 0.0571715145436308
$ws.Range("F25").Value = 39.461538461538403
$ws.Range("D26").Value = 0.15
$ws.Range("E26").Value = 0.25
$ws.Range("F26").Value = 39.461538461538403
$ws.Range("D27").Value = 58.290155440414502
$ws.Range("E27").Value = 0.0571715145436308
$ws.Range("F27").Value = 39.461538461538403
$ws.Range("D28").Value = 1.5385
$ws.Range("E28").Value = 0.252584358966155
$ws.Range("F28").Value = 6.31460897415388

# --- Row 2: explicit (non-shared) formulas ---
$ws.Range("G2").Formula = "=E2*D2"
$ws.Range("H2").Formula = "=F2*D2"

# --- Rows 3-28: shared formulas (G3 / H3 become the shared-formula masters) ---
$ws.Range("G3:G28").Formula = "=E3*D3"
$ws.Range("H3:H28").Formula = "=F3*D3"

# --- Selection / view state (best effort) ---
$ws.Range("G2:H28").Select()
